$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2  = @(15.2503422516503, 0.14960932783464201, 15.1056375175573)
    3  = @(14.910502642942401, 0.091427814115147901, 12.026642027228499)
    4  = @(15.070097494195799, 0.128488329257539, 13.9594652099748)
    5  = @(15.2325937252452, 0.169644837392775, 16.0985387322614)
    6  = @(15.207310177674101, 0.18193010357646799, 16.7001602676605)
    7  = @(15.4512604760295, 0.16832949989302001, 16.142105232665799)
    8  = @(15.44384358271, 0.17274330060533999, 16.359086821622)
    9  = @(14.8515294961508, 0.20056219205713199, 17.453874350932001)
    10 = @(15.2218146262712, 0.24931600820071301, 20.0767077231713)
    11 = @(15.0731759385783, 0.1914247251193, 17.107824225254099)
    12 = @(15.580306614356999, 0.22327100457638999, 18.953703535997999)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
}

$ws.Range("L18").Select()
